$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'94.084.14"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.45%  '
$ws.Range('D3').Value = "'3.071.49"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'235.67"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = "'609.54"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E7').Value = '  +1.89%  '
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').Value = "'0.999"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').Value = "'0.810"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.86%  '
$ws.Range('D11').Value = "'3.071.79"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').Value = "'93.995.10"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('E14').Value = '  -2.40%  '
$ws.Range('D15').Value = "'33.84"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = "'3.652.58"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').Value = "'3.069.28"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').Value = "'3.56"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Value = "'14.34"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').Value = "'5.73"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').Value = "'444.08"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').Value = "'8.81"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.72%  '
$ws.Range('E24').Value = '  -2.67%  '
$ws.Range('D25').Value = "'8.34"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.57%  '
$ws.Range('D26').Value = "'5.50"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.82%  '
$ws.Range('D27').Value = "'84.54"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').Value = "'11.92"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.16%  '
$ws.Range('D29').Value = "'3.242.16"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').Value = "'0.249"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.61%  '
$ws.Range('E32').Value = '  +7.31%  '
$ws.Range('E33').Value = '  -6.02%  '
$ws.Range('D34').Value = "'1.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.75%  '
$ws.Range('D35').Value = "'8.93"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').Value = "'7.54"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.86%  '
$ws.Range('D37').Value = "'25.42"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('E38').Value = '  -4.44%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').Value = "'479.30"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').Value = "'24.04"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.73%  '
$ws.Range('E42').Value = '  +1.83%  '
$ws.Range('D43').Value = "'3.75"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('D44').Value = "'1.25"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.79%  '
$ws.Range('D46').Value = "'3.08"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.96%  '
$ws.Range('D47').Value = "'161.56"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('D50').Value = "'43.62"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('E51').Value = '  +0.12%  '
